$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 413.9565
$ws.Range("I33").Value = 332.77274
$ws.Range("K33").Value = 332.77274
$ws.Range("M33").Value = -103.77274

$ws.Range("H40").Value = 371940.8
$ws.Range("I40").Value = 835723.5600000001
$ws.Range("J40").Value = 914.6
$ws.Range("K40").Value = 835723.5600000001
$ws.Range("L40").Value = 914.6
$ws.Range("M40").Value = -835548.5600000001
$ws.Range("N40").Value = -1264.6

$ws.Range("H51").Value = 4999.5625
$ws.Range("J51").Value = 4999.5625
$ws.Range("L51").Value = 4999.5625
$ws.Range("N51").Value = -5967.5625

$ws.Range("H86").Value = 15876721
$ws.Range("I86").Value = 3112.75
$ws.Range("J86").Value = 25645096
$ws.Range("K86").Value = 3112.75
$ws.Range("L86").Value = 25645096
$ws.Range("M86").Value = -1989.75
$ws.Range("N86").Value = -25647342

$ws.Range("H89").Value = 15876721
$ws.Range("I89").Value = 3112.75
$ws.Range("J89").Value = 25645096
$ws.Range("K89").Value = 15563.75
$ws.Range("L89").Value = 128225480
$ws.Range("M89").Value = -9947.75
$ws.Range("N89").Value = -128236712

$ws.Range("H98").Value = 3240.4167
$ws.Range("I98").Value = 3307.7273
$ws.Range("J98").Value = 2500
$ws.Range("K98").Value = 3307.7273
$ws.Range("L98").Value = 2500
$ws.Range("M98").Value = -1809.7273
$ws.Range("N98").Value = -5496

$ws.Range("H107").Value = 2017
$ws.Range("I107").Value = 687.5
$ws.Range("K107").Value = 687.5
$ws.Range("M107").Value = 1232.5

$ws.Range("H113").Value = 9174.333000000001
$ws.Range("I113").Value = 9174.333000000001
$ws.Range("K113").Value = 9174.333000000001
$ws.Range("M113").Value = -5920.333000000001

$ws.Range("H122").Value = 3240.4167
$ws.Range("I122").Value = 3307.7273
$ws.Range("J122").Value = 2500
$ws.Range("K122").Value = 9923.1819
$ws.Range("L122").Value = 7500
$ws.Range("M122").Value = -7473.1819
$ws.Range("N122").Value = -12400

$ws.Range("H132").Value = 3382.6191
$ws.Range("I132").Value = 3476.6858
$ws.Range("K132").Value = 10430.0574
$ws.Range("M132").Value = -7900.057400000002

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H110").Value = 17108.8
$ws.Range("I110").Value = 21551
$ws.Range("J110").Value = 3041.8333
$ws.Range("K110").Value = 21551
$ws.Range("L110").Value = 3041.8333
$ws.Range("M110").Value = -19506
$ws.Range("N110").Value = -7131.8333

$ws.Range("H122").Value = 7580297
$ws.Range("I122").Value = 2585.75
$ws.Range("J122").Value = 9264233
$ws.Range("K122").Value = 7757.25
$ws.Range("L122").Value = 27792699
$ws.Range("M122").Value = -5307.25
$ws.Range("N122").Value = -27797599

$ws.Range("H132").Value = 31259432
$ws.Range("I132").Value = 5394.3
$ws.Range("J132").Value = 500070000
$ws.Range("K132").Value = 16182.9
$ws.Range("L132").Value = 1500210000
$ws.Range("M132").Value = -13652.9
$ws.Range("N132").Value = -1500215060

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1511.3529
$ws.Range("I94").Value = 1249.5834
$ws.Range("J94").Value = 2139.6
$ws.Range("K94").Value = 1249.5834
$ws.Range("L94").Value = 2139.6
$ws.Range("M94").Value = -798.5834
$ws.Range("N94").Value = -3041.6

$ws.Range("H99").Value = 2568.1428
$ws.Range("I99").Value = 2568.1428
$ws.Range("K99").Value = 2568.1428
$ws.Range("M99").Value = -1070.1428

$ws.Range("H140").Value = 98092.14
$ws.Range("J140").Value = 98092.14
$ws.Range("L140").Value = 98092.14
$ws.Range("N140").Value = -108452.14

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 10635.3
$ws.Range("I22").Value = 17071.5
$ws.Range("K22").Value = 17071.5
$ws.Range("M22").Value = -16721.5

$ws.Range("H31").Value = 40325212
$ws.Range("J31").Value = 125004350
$ws.Range("L31").Value = 125004350
$ws.Range("N31").Value = -125004940

$ws.Range("H34").Value = 40325212
$ws.Range("J34").Value = 125004350
$ws.Range("L34").Value = 125004350
$ws.Range("N34").Value = -125004754

$ws.Range("H132").Value = 2091.3333
$ws.Range("I132").Value = 2091.3333
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 6273.999899999999
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = $null
$ws.Range("N132").Value = -3743.999899999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 2296.3333
$ws.Range("J5").Value = 2770.9412
$ws.Range("L5").Value = 8312.8236
$ws.Range("N5").Value = -8536.8236

$ws.Range("H55").Value = 902920
$ws.Range("I55").Value = 1501033.4
$ws.Range("K55").Value = 4503100.199999999
$ws.Range("M55").Value = -4502923.199999999

$ws.Range("H107").Value = 833.44446
$ws.Range("J107").Value = 960.6
$ws.Range("L107").Value = 2881.8
$ws.Range("N107").Value = -6721.8

$ws.Range("H113").Value = 1065.5834
$ws.Range("I113").Value = 550
$ws.Range("J113").Value = 1237.4445
$ws.Range("K113").Value = 1650
$ws.Range("L113").Value = 3712.3335
$ws.Range("M113").Value = 520
$ws.Range("N113").Value = -8052.333500000001

$ws.Range("H135").Value = 2296.3333
$ws.Range("J135").Value = 2770.9412
$ws.Range("L135").Value = 24938.4708
$ws.Range("N135").Value = -30008.4708

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H57").Value = 41250
$ws.Range("I57").Value = 15000
$ws.Range("J57").Value = 50000
$ws.Range("K57").Value = 15000
$ws.Range("L57").Value = 50000
$ws.Range("M57").Value = -14180
$ws.Range("N57").Value = -51640

$ws.Range("H88").Value = 67436.60000000001
$ws.Range("I88").Value = 64990
$ws.Range("K88").Value = 64990
$ws.Range("M88").Value = -64539

$ws.Range("H91").Value = 67436.60000000001
$ws.Range("I91").Value = 64990
$ws.Range("K91").Value = 64990
$ws.Range("M91").Value = -63430

$ws.Range("H99").Value = 27350.7
$ws.Range("I99").Value = 13917.833
$ws.Range("K99").Value = 13917.833
$ws.Range("M99").Value = -11671.833

$ws.Range("H132").Value = 13441.541
$ws.Range("I132").Value = 11559.29
$ws.Range("J132").Value = 23166.5
$ws.Range("K132").Value = 34677.87
$ws.Range("L132").Value = 69499.5
$ws.Range("M132").Value = -32147.87
$ws.Range("N132").Value = -74559.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1254.6552
$ws.Range("I16").Value = 1299.6296
$ws.Range("J16").Value = 647.5
$ws.Range("K16").Value = 1299.6296
$ws.Range("L16").Value = 647.5
$ws.Range("M16").Value = -1129.6296
$ws.Range("N16").Value = -987.5

$ws.Range("H40").Value = 6327.3125
$ws.Range("I40").Value = 6148.6924
$ws.Range("K40").Value = 6148.6924
$ws.Range("M40").Value = -6012.6924

$ws.Range("H46").Value = 1546.3889
$ws.Range("J46").Value = 3857.1428
$ws.Range("L46").Value = 3857.1428
$ws.Range("N46").Value = -4233.1428

$ws.Range("H55").Value = 423.83334
$ws.Range("I55").Value = 413.85715
$ws.Range("J55").Value = 437.8
$ws.Range("K55").Value = 413.85715
$ws.Range("L55").Value = 437.8
$ws.Range("M55").Value = -240.85715
$ws.Range("N55").Value = -783.8

$ws.Range("H100").Value = 2212.2903
$ws.Range("I100").Value = 1951.1724
$ws.Range("K100").Value = 1951.1724
$ws.Range("M100").Value = -1410.1724

$ws.Range("H132").Value = 55558556
$ws.Range("I132").Value = 3178.7058
$ws.Range("K132").Value = 9536.117400000001
$ws.Range("M132").Value = -7006.117400000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H14").Value = 647
$ws.Range("I14").Value = 647
$ws.Range("K14").Value = 647
$ws.Range("M14").Value = -479

$ws.Range("H54").Value = 38495
$ws.Range("J54").Value = 38495
$ws.Range("L54").Value = 38495
$ws.Range("N54").Value = -39535

$ws.Range("H81").Value = 3780.0435
$ws.Range("I81").Value = 3282.9048
$ws.Range("J81").Value = 9000
$ws.Range("K81").Value = 6565.8096
$ws.Range("L81").Value = 18000
$ws.Range("M81").Value = -5504.8096
$ws.Range("N81").Value = -20122

$ws.Range("H84").Value = 3780.0435
$ws.Range("I84").Value = 3282.9048
$ws.Range("J84").Value = 9000
$ws.Range("K84").Value = 32829.048
$ws.Range("L84").Value = 90000
$ws.Range("M84").Value = -27525.048
$ws.Range("N84").Value = -100608

$ws.Range("H113").Value = 1057.2142
$ws.Range("I113").Value = 750.0833
$ws.Range("J113").Value = 2900
$ws.Range("K113").Value = 2250.2499
$ws.Range("L113").Value = 8700
$ws.Range("M113").Value = -80.2498999999998
$ws.Range("N113").Value = -13040

$ws.Range("H122").Value = 5558917
$ws.Range("I122").Value = 3110.1924
$ws.Range("K122").Value = 9330.5772
$ws.Range("M122").Value = -6880.5772

$ws.Range("H126").Value = 2519.7742
$ws.Range("I126").Value = 1793.7727
$ws.Range("J126").Value = 4294.4443
$ws.Range("K126").Value = 5381.3181
$ws.Range("L126").Value = 12883.3329
$ws.Range("M126").Value = -2911.3181
$ws.Range("N126").Value = -17823.3329
